$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Status column for several in-flight LLM tasks to "Completed" ---
$ws.Range("C7").Value  = "Completed"
$ws.Range("C8").Value  = "Completed"
$ws.Range("C9").Value  = "Completed"
$ws.Range("C10").Value = "Completed"
$ws.Range("C13").Value = "Completed"
$ws.Range("C14").Value = "Completed"

# --- Row 9's start date moved earlier ---
$ws.Range("D9").Value = 45600

# --- Fill in the previously-blank row 15 (new task row) ---
$ws.Range("A15").Value = "Viết báo cáo cuối kì"
$ws.Range("B15").Value = "Tâm Tăng Thiện Bảo Nguyễn Lê Gia"
$ws.Range("C15").Value = "In progress"

# Owner cell on the new row gets the same mailto hyperlink style used
# throughout the rest of the "Owner" column.
$ws.Hyperlinks.Add($ws.Range("B15"), "mailto:21521408@gm.uit.edu.vn")
$ws.Range("B15").Font.Underline = $true
$ws.Range("B15").Font.ThemeColor = 1
